$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 11:
# value 45174 (2023-09-05) -> 45175 (2023-09-06)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
